# Update column F (dSF) values for specific rows as per repulled data / mean calculation fix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1
$ws.Range("F6").Value = -2
$ws.Range("F7").Value = -5
$ws.Range("F12").Value = 0
$ws.Range("F22").Value = 1
$ws.Range("F24").Value = -2
$ws.Range("F28").Value = -2
$ws.Range("F30").Value = 1
$ws.Range("F31").Value = 4
$ws.Range("F32").Value = 1
$ws.Range("F39").Value = 2
$ws.Range("F47").Value = 0
